$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2033.75
$ws.Cells.Item(40, 9).Value = 2294.1765
$ws.Cells.Item(40, 10).Value = 1401.2858
$ws.Cells.Item(40, 11).Value = 2294.1765
$ws.Cells.Item(40, 12).Value = 1401.2858
$ws.Cells.Item(40, 13).Value = -2119.1765
$ws.Cells.Item(40, 14).Value = -1751.2858
$ws.Cells.Item(137, 8).Value = 1738.4916
$ws.Cells.Item(137, 9).Value = 1256.7872
$ws.Cells.Item(137, 11).Value = 3770.3616
$ws.Cells.Item(137, 13).Value = -1220.3616
$ws.Cells.Item(138, 8).Value = 3344.7058
$ws.Cells.Item(138, 9).Value = 1542.0312
$ws.Cells.Item(138, 11).Value = 4626.0936
$ws.Cells.Item(138, 13).Value = 513.9063999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 8009.28
$ws.Cells.Item(61, 9).Value = 4590.846
$ws.Cells.Item(61, 10).Value = 20129.182
$ws.Cells.Item(61, 11).Value = 4590.846
$ws.Cells.Item(61, 12).Value = 20129.182
$ws.Cells.Item(61, 13).Value = -4378.846
$ws.Cells.Item(61, 14).Value = -20553.182
$ws.Cells.Item(74, 8).Value = 7835.2285
$ws.Cells.Item(74, 9).Value = 5921.6665
$ws.Cells.Item(74, 10).Value = 14293.5
$ws.Cells.Item(74, 11).Value = 5921.6665
$ws.Cells.Item(74, 12).Value = 14293.5
$ws.Cells.Item(74, 13).Value = -5047.6665
$ws.Cells.Item(74, 14).Value = -16041.5
$ws.Cells.Item(77, 8).Value = 7835.2285
$ws.Cells.Item(77, 9).Value = 5921.6665
$ws.Cells.Item(77, 10).Value = 14293.5
$ws.Cells.Item(77, 11).Value = 29608.3325
$ws.Cells.Item(77, 12).Value = 71467.5
$ws.Cells.Item(77, 13).Value = -25240.3325
$ws.Cells.Item(77, 14).Value = -80203.5
$ws.Cells.Item(102, 8).Value = 1483321.4
$ws.Cells.Item(102, 9).Value = 1765020.6
$ws.Cells.Item(102, 11).Value = 1765020.6
$ws.Cells.Item(102, 13).Value = -1763398.6
$ws.Cells.Item(122, 8).Value = 2907762
$ws.Cells.Item(122, 9).Value = 634.1081
$ws.Cells.Item(122, 10).Value = 20835050
$ws.Cells.Item(122, 11).Value = 1902.3243
$ws.Cells.Item(122, 12).Value = 62505150
$ws.Cells.Item(122, 13).Value = 547.6756999999998
$ws.Cells.Item(122, 14).Value = -62510050
$ws.Cells.Item(136, 8).Value = 8009.28
$ws.Cells.Item(136, 9).Value = 4590.846
$ws.Cells.Item(136, 10).Value = 20129.182
$ws.Cells.Item(136, 11).Value = 13772.538
$ws.Cells.Item(136, 12).Value = 60387.546
$ws.Cells.Item(136, 13).Value = -11222.538
$ws.Cells.Item(136, 14).Value = -65487.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1599.8334
$ws.Cells.Item(94, 9).Value = 1544.1428
$ws.Cells.Item(94, 10).Value = 1677.8
$ws.Cells.Item(94, 11).Value = 1544.1428
$ws.Cells.Item(94, 12).Value = 1677.8
$ws.Cells.Item(94, 13).Value = -1093.1428
$ws.Cells.Item(94, 14).Value = -2579.8
$ws.Cells.Item(105, 8).Value = 1841709.5
$ws.Cells.Item(105, 9).Value = 2086523.4
$ws.Cells.Item(105, 10).Value = 5605.5
$ws.Cells.Item(105, 11).Value = 2086523.4
$ws.Cells.Item(105, 12).Value = 5605.5
$ws.Cells.Item(105, 13).Value = -2084776.4
$ws.Cells.Item(105, 14).Value = -9099.5
$ws.Cells.Item(134, 8).Value = 17662.604
$ws.Cells.Item(134, 9).Value = 1622.6666
$ws.Cells.Item(134, 10).Value = 68990.39999999999
$ws.Cells.Item(134, 11).Value = 4867.9998
$ws.Cells.Item(134, 12).Value = 206971.2
$ws.Cells.Item(134, 13).Value = -2332.9998
$ws.Cells.Item(134, 14).Value = -212041.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(43, 8).Value = 42500
$ws.Cells.Item(43, 10).Value = 42500
$ws.Cells.Item(43, 12).Value = 42500
$ws.Cells.Item(43, 14).Value = -42868
$ws.Cells.Item(101, 8).Value = 42500
$ws.Cells.Item(101, 10).Value = 42500
$ws.Cells.Item(101, 12).Value = 42500
$ws.Cells.Item(101, 14).Value = -48990
$ws.Cells.Item(111, 8).Value = 53266.668
$ws.Cells.Item(111, 10).Value = 53266.668
$ws.Cells.Item(111, 12).Value = 53266.668
$ws.Cells.Item(111, 14).Value = -61446.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 2717.5
$ws.Cells.Item(46, 10).Value = 2717.5
$ws.Cells.Item(46, 12).Value = 8152.5
$ws.Cells.Item(46, 14).Value = -8334.5
$ws.Cells.Item(86, 8).Value = 718.25
$ws.Cells.Item(86, 9).Value = 691.3333
$ws.Cells.Item(86, 10).Value = 799
$ws.Cells.Item(86, 11).Value = 2073.9999
$ws.Cells.Item(86, 12).Value = 2397
$ws.Cells.Item(86, 13).Value = -887.9998999999998
$ws.Cells.Item(86, 14).Value = -4769
$ws.Cells.Item(89, 8).Value = 718.25
$ws.Cells.Item(89, 9).Value = 691.3333
$ws.Cells.Item(89, 10).Value = 799
$ws.Cells.Item(89, 11).Value = 6221.9997
$ws.Cells.Item(89, 12).Value = 7191
$ws.Cells.Item(89, 13).Value = -293.9997000000003
$ws.Cells.Item(89, 14).Value = -19047
$ws.Cells.Item(131, 8).Value = 17452.293
$ws.Cells.Item(131, 9).Value = 963.7619
$ws.Cells.Item(131, 10).Value = 26810.648
$ws.Cells.Item(131, 11).Value = 2891.2857
$ws.Cells.Item(131, 12).Value = 80431.944
$ws.Cells.Item(131, 13).Value = 2148.7143
$ws.Cells.Item(131, 14).Value = -90511.944

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(112, 8).Value = 48230
$ws.Cells.Item(112, 10).Value = 48230
$ws.Cells.Item(112, 12).Value = 48230
$ws.Cells.Item(112, 14).Value = -50446
$ws.Cells.Item(126, 8).Value = 2463.4
$ws.Cells.Item(126, 9).Value = 1610.1818
$ws.Cells.Item(126, 11).Value = 4830.5454
$ws.Cells.Item(126, 13).Value = -2360.5454
$ws.Cells.Item(134, 8).Value = 42177.816
$ws.Cells.Item(134, 10).Value = 42177.816
$ws.Cells.Item(134, 12).Value = 126533.448
$ws.Cells.Item(134, 14).Value = -131603.448
$ws.Cells.Item(135, 8).Value = 53813.332
$ws.Cells.Item(135, 10).Value = 53813.332
$ws.Cells.Item(135, 12).Value = 53813.332
$ws.Cells.Item(135, 14).Value = -63953.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(95, 8).Value = 26108
$ws.Cells.Item(95, 10).Value = 26108
$ws.Cells.Item(95, 12).Value = 26108
$ws.Cells.Item(95, 14).Value = -31600
$ws.Cells.Item(106, 8).Value = 97066
$ws.Cells.Item(106, 10).Value = 97066
$ws.Cells.Item(106, 12).Value = 97066
$ws.Cells.Item(106, 14).Value = -99590
$ws.Cells.Item(108, 8).Value = 49900
$ws.Cells.Item(108, 10).Value = 49900
$ws.Cells.Item(108, 12).Value = 49900
$ws.Cells.Item(108, 14).Value = -57580
$ws.Cells.Item(132, 8).Value = 3980.6416
$ws.Cells.Item(132, 9).Value = 4137.5
$ws.Cells.Item(132, 10).Value = 3498
$ws.Cells.Item(132, 11).Value = 12412.5
$ws.Cells.Item(132, 12).Value = 10494
$ws.Cells.Item(132, 13).Value = -9882.5
$ws.Cells.Item(132, 14).Value = -15554
$ws.Cells.Item(136, 8).Value = 3720.9692
$ws.Cells.Item(136, 9).Value = 2143.4348
$ws.Cells.Item(136, 10).Value = 7540.263
$ws.Cells.Item(136, 11).Value = 6430.3044
$ws.Cells.Item(136, 12).Value = 22620.789
$ws.Cells.Item(136, 13).Value = -3880.3044
$ws.Cells.Item(136, 14).Value = -27720.789
$ws.Cells.Item(139, 8).Value = 60715
$ws.Cells.Item(139, 10).Value = 60715
$ws.Cells.Item(139, 12).Value = 60715
$ws.Cells.Item(139, 14).Value = -70995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 25083.111
$ws.Cells.Item(123, 10).Value = 25083.111
$ws.Cells.Item(123, 12).Value = 25083.111
$ws.Cells.Item(123, 14).Value = -34883.111
$ws.Cells.Item(125, 8).Value = 64857.5
$ws.Cells.Item(125, 10).Value = 64857.5
$ws.Cells.Item(125, 12).Value = 64857.5
$ws.Cells.Item(125, 14).Value = -74697.5
$ws.Cells.Item(136, 8).Value = 2904.2908
$ws.Cells.Item(136, 9).Value = 1165.0536
$ws.Cells.Item(136, 10).Value = 6150.8667
$ws.Cells.Item(136, 11).Value = 3495.1608
$ws.Cells.Item(136, 12).Value = 18452.6001
$ws.Cells.Item(136, 13).Value = -945.1607999999997
$ws.Cells.Item(136, 14).Value = -23552.6001

Write-Output "Updated 173 cells across 8 sheets"